$wb = $excel.ActiveWorkbook

# --- Sheet "commA-v2" (additional data examples) ---
$ws = $wb.Worksheets.Item("commA-v2")
$ws.Range("B1").Value = 12
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 3
$ws.Range("A7").Value = 2
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = 4
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = 5
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = 5
$ws.Range("A16").Value = 8
$ws.Range("B16").Value = 6
$ws.Range("A17").Value = 10
$ws.Range("B17").Value = 12
$ws.Activate()
$ws.Range("A18").Select()

# --- Sheet "commB-v2" (additional data examples) ---
$ws2 = $wb.Worksheets.Item("commB-v2")
$ws2.Range("B1").Value = 12
$ws2.Range("A5").Value = 1
$ws2.Range("B5").Value = 5
$ws2.Range("B6").Value = 3
$ws2.Range("A7").Value = 2
$ws2.Range("A8").Value = 3
$ws2.Range("B8").Value = 4
$ws2.Range("A9").Value = 3
$ws2.Range("B9").Value = 5
$ws2.Range("A10").Value = 4
$ws2.Range("B10").Value = 5
$ws2.Range("A16").Value = 8
$ws2.Range("B16").Value = 6
$ws2.Range("A17").Value = 10
$ws2.Range("B17").Value = 12
$ws2.Range("A18").Value = 9
$ws2.Range("B18").Value = 11

# "commB-v2" is the active/selected tab in the final workbook state.
$ws2.Activate()
$ws2.Range("F22").Select()
